$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new row 2 ("Player wants to log in to their account" use case),
# shifting all the existing rows (old row 2 -> row 3, old row 3 -> row 4, ...)
# down by one. Copy the format from the row that ends up directly below the
# new row (row 3, which carries the standard s="1"/s="4" style pattern used
# by every other data row) so the new row picks up identical styling.
# ---------------------------------------------------------------------------
$ws.Rows("2:2").Insert()
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new row 2.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "Player wants to log in to their account"
$ws.Range("C2").Value = "Username and Password"
$ws.Range("D2").Value = "Player"
$ws.Range("E2").Value = "Player is taken to their dashboard"
$ws.Range("F2").Value = "Player"

# The two rows that follow the new row (now rows 4 and 5, formerly rows 3
# and 4) gain a more specific "Actor Receiving Output" label. (F5 is set
# before F4 so the shared-string table picks up the two new strings in the
# same order as the reference workbook.)
$ws.Range("F5").Value = "Player / Company that the User's bank account is from"
$ws.Range("F4").Value = "Player / SMS Respondent"

# ---------------------------------------------------------------------------
# Row heights (re-autofit by Excel after the content/row changes).
# ---------------------------------------------------------------------------
$ws.Rows("2:2").RowHeight = 29.45
$ws.Rows("3:3").RowHeight = 29.45
$ws.Rows("4:4").RowHeight = 31.15
$ws.Rows("5:5").RowHeight = 45
$ws.Rows("6:6").RowHeight = 71.25
$ws.Rows("7:7").RowHeight = 31.15
$ws.Rows("8:8").RowHeight = 36.6
$ws.Rows("9:9").RowHeight = 35.45
$ws.Rows("10:10").RowHeight = 33.6
$ws.Rows("11:11").RowHeight = 33.6
$ws.Rows("12:12").RowHeight = 33.6
$ws.Rows("13:13").RowHeight = 61.15
$ws.Rows("14:14").RowHeight = 45

# Column D widened to fit the longer "Actor" values.
$ws.Columns("D:D").ColumnWidth = 14.75

# ---------------------------------------------------------------------------
# View state: selection moved to F5, no frozen/scrolled top-left cell.
# ---------------------------------------------------------------------------
[void]$ws.Range("F5").Select()

Write-Output "done"
